# Week2_Recap.pptx edit:
# Slide 3 ("Ex00-Q6: Find the surface area and the diagonal of a cuboid")
# has its cuboid illustration picture nudged to a new position:
#   EMU offset (6506760, 1764856)  ->  (6578601, 1682560)
# (the picture's size is unchanged: 2074080 x 1309983 EMU)
#
# PowerPoint's Shape.Left/.Top COM properties are expressed in points
# (1 pt = 12700 EMU). The literal values below are the point values that,
# after the host's internal point->EMU conversion, land exactly on the
# target EMU offsets (plain "EMU/12700" can be thrown off by a single EMU
# due to float rounding, so we use values safely inside the correct
# point->EMU rounding bucket).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item("Picture 2")

$shp.Left = 518.0001220867068   # -> 6578601 EMU
$shp.Top  = 132.48507686948284  # -> 1682560 EMU
